# Weekly update for "Fruta, Terminal La Palmera de La Serena - Uva".
#
# A new weekly price record is inserted as row 121 (pushing the existing
# rows 121-146 down to 122-147), matching the author's "Fruta / hortaliza,
# semanal" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 121; everything below shifts down one row.
$ws.Rows.Item(121).Insert()

# Populate the new row with the new weekly record.
$ws.Cells.Item(121, 1).Value2  = 8
$ws.Cells.Item(121, 2).Value2  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(121, 3).Value2  = "Coquimbo"
$ws.Cells.Item(121, 4).Value2  = 44985
$ws.Cells.Item(121, 5).Value2  = 4
$ws.Cells.Item(121, 6).Value2  = "Fruta"
$ws.Cells.Item(121, 7).Value2  = 100109
$ws.Cells.Item(121, 8).Value2  = "Uva"
$ws.Cells.Item(121, 9).Value2  = 100109001
$ws.Cells.Item(121, 10).Value2 = "Uva"
$ws.Cells.Item(121, 11).Value2 = "Thompson seedless"
$ws.Cells.Item(121, 12).Value2 = "Primera"
$ws.Cells.Item(121, 13).Value2 = 540
$ws.Cells.Item(121, 14).Value2 = 9500
$ws.Cells.Item(121, 15).Value2 = 10000
$ws.Cells.Item(121, 16).Value2 = 9750
$ws.Cells.Item(121, 17).Value2 = "`$/bandeja 18 kilos"
$ws.Cells.Item(121, 18).Value2 = "Provincia del Elquí"
$ws.Cells.Item(121, 19).Value2 = 542
$ws.Cells.Item(121, 20).Value2 = 18
